$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-51 down to 43-52
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record
$ws.Cells.Item(42, 1).Value = 5
$ws.Cells.Item(42, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(42, 3).Value = "Maule"
$ws.Cells.Item(42, 4).Value = 44798
$ws.Cells.Item(42, 5).Value = 7
$ws.Cells.Item(42, 6).Value = 100112040
$ws.Cells.Item(42, 7).Value = "Cilantro"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 150
$ws.Cells.Item(42, 11).Value = 13000
$ws.Cells.Item(42, 12).Value = 13000
$ws.Cells.Item(42, 13).Value = 13000
$ws.Cells.Item(42, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(42, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(42, 16).Value = 361
$ws.Cells.Item(42, 17).Value = 36
$ws.Cells.Item(42, 18).Value = "Hortaliza"

# Ensure the Date column (D) keeps the date number format for the new row
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
